$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "key"
$ws.Range("B4").Font.Bold = $true

$ws.Range("B5").Value = "Value"
$ws.Range("B5").Font.Bold = $false
